$d = $word.ActiveDocument

$d.Content.Find.Execute("05-07-16", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "23-02-17", 2)

$d.Content.Find.Execute("SV1606220059", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "SV1702200002", 2)

$d.Content.Find.Execute("PAK'nSAVE Mill Street", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "New World Rototuna", 2)

$d.Content.Find.Execute("G14M91831", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "AG22016765E0", 2)

$d.Content.Find.Execute("110467", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "900015159", 2)

$d.Content.Find.Execute("MJ0049368", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "AK25005042E0", 2)

$d.Content.Find.Execute("900004369", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "900017086", 2)

$d.Content.Find.Execute("failed calibration. Replaced and tested weight and scanning items. Working fine. Told staff to get calibrated.", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "printer destroyed multiple ink ribbons this morning. Replaced with new printer, configured and tested. Working fine.", 2)
